# Project_LibraryAnalysis_TaskManager.xlsx
# Commit: "Add my Python Pandas out of my Jupyter notebook."
#
# Content change: the "Cleaning Data" task list gains a new subtask row
# ("Remove duplicates.") inserted right after "Limit data set to ... years.",
# and that first subtask's wording is tightened from 10 years to 5 years.
# Every row below shifts down by one to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tighten the retention window mentioned in the first "Cleaning Data" subtask.
$ws.Range("B2").Value = "Limit data set to 5 years."

# Insert a new row right below it for the additional cleaning subtask; this
# pushes every following row (Data Analysis tasks, blank rows, validation
# rows) down by one, which is exactly what the target sheet shows.
$ws.Range("A3").EntireRow.Insert()
$ws.Range("A3").Value = "Cleaning Data"
$ws.Range("B3").Value = "Remove duplicates."

# Restore the active selection to A2 (top of the frozen pane body) like the
# saved workbook shows.
$ws.Range("A2").Select() | Out-Null
